$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvestorKyc")

# Add new "Agreement Committed Amount" column header in R1
$ws.Range("R1").Value = "Agreement Committed Amount"

# Rename "Full Name" header (column B) to "Investing Entity"
$ws.Range("B1").Value = "Investing Entity"

# Fill in Agreement Committed Amount values for the 4 data rows
$ws.Range("R2").Value = 1000000
$ws.Range("R3").Value = 2000000
$ws.Range("R4").Value = 3000000
$ws.Range("R5").Value = 4000000

# Update the active selection to R6
$ws.Range("R6").Select()
